$wb = $excel.ActiveWorkbook

# --- Sheet "Metadata": refresh the generation timestamp ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-10-30T16:36:55+00:00"

# --- Sheet "Elements": widen column K ("Type(s)") to fit the new URL content ---
$els = $wb.Worksheets.Item("Elements")
$els.Columns.Item(11).ColumnWidth = 55.52

# --- Sheet "Elements": add row 8 for "AutorisationExercice.Professionnel" ---
# Pre-seed the numeric-looking cardinality cells (Min/Max columns) with a
# leading apostrophe so Excel stores them as TEXT ("0"/"1"), matching how the
# rest of the sheet stores these values (as shared strings, not numbers).
$els.Range("F8").Value = "'0"
$els.Range("G8").Value = "'1"
$els.Range("AG8").Value = "'0"
$els.Range("AH8").Value = "'1"

# Copy row 7's formatting onto row 8: this materialises every cell A8:AJ8
# (including the ones that stay empty) with the same style used by every
# other data row.
$els.Range("A7:AJ7").Copy()
$els.Range("A8:AJ8").PasteSpecial(-4122)

# Fill in the row's textual content.
$els.Range("A8").Value = "AutorisationExercice.Professionnel"
$els.Range("B8").Value = "AutorisationExercice.Professionnel"
$els.Range("K8").Value = "https://interop.esante.gouv.fr/ig/mos/StructureDefinition/Professionnel`n"
$els.Range("L8").Value = "Lien vers la classe Professionnel"
$els.Range("M8").Value = "Lien vers la classe Professionnel"
$els.Range("AF8").Value = "AutorisationExercice.Professionnel"

# These three columns (Binding Strength / Description / Value Set) are not
# populated for this element (unlike the sibling Coding rows above it), so
# clear out whatever row 7's formatting copy carried over.
$els.Range("X8").Value = ""
$els.Range("Y8").Value = ""
$els.Range("Z8").Value = ""
